$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the usage date / usage count header field names to make them
# consistent ("R4_Month" -> "Usage_Date", "R4_Count" -> "Usage_Count").
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Reflect the new selection left behind in the saved view state.
$ws.Range("K1:L1").Select()
